# Adds 9 new "Function" dictionary entries (and their blank spacer rows)
# to the bottom of the Data Dictionary sheet, mirroring the existing
# A/B (name, merged) | C/D (description, merged) | E (category, merged) |
# F/G (name repeated, merged) row-pair layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert 9 fresh two-row blocks (copied from the clean row 4:5
#    template, which already carries the correct borders/alignment
#    and no stray custom height) right after the current last row.
# ---------------------------------------------------------------
$ws.Rows("4:5").Copy()
$ws.Rows("212:213").Insert()
$ws.Rows("4:5").Copy()
$ws.Rows("214:215").Insert()
$ws.Rows("4:5").Copy()
$ws.Rows("216:217").Insert()
$ws.Rows("4:5").Copy()
$ws.Rows("218:219").Insert()
$ws.Rows("4:5").Copy()
$ws.Rows("220:221").Insert()
$ws.Rows("4:5").Copy()
$ws.Rows("222:223").Insert()
$ws.Rows("4:5").Copy()
$ws.Rows("224:225").Insert()
$ws.Rows("4:5").Copy()
$ws.Rows("226:227").Insert()
$ws.Rows("4:5").Copy()
$ws.Rows("228:229").Insert()

# ---------------------------------------------------------------
# 2. Re-create the A:B / C:D / E:E / F:G merges on each new pair
#    (Insert-from-copy does not carry merged-cell ranges over).
# ---------------------------------------------------------------
$pairs = @(212, 214, 216, 218, 220, 222, 224, 226, 228)
foreach ($r in $pairs) {
    $r2 = $r + 1
    $ws.Range("A${r}:B${r2}").Merge()
    $ws.Range("C${r}:D${r2}").Merge()
    $ws.Range("E${r}:E${r2}").Merge()
    $ws.Range("F${r}:G${r}").Merge()
}

# ---------------------------------------------------------------
# 3. Fill in the values, matching the original authoring order so
#    the shared-string table is rebuilt in the same sequence:
#    names first (rows 212/214/216/218), then their descriptions,
#    then the remaining names (222 before 220!), then the rest of
#    the descriptions in row order.
# ---------------------------------------------------------------
$ws.Range("A212").Value = "allItemsTotal"
$ws.Range("A214").Value = "allLendedItemsTotal"
$ws.Range("A216").Value = "allItems"
$ws.Range("A218").Value = "allLendedItems"

$ws.Range("C212").Value = "Consulta que devuelve la cantidad de items en la tabla Items"
$ws.Range("C214").Value = "Consulta que devuelve la cantidad de items en la tabla de PersonLendItem"
$ws.Range("C216").Value = "Consulta que devuelve todos los items en orden alfabetico, filtrable por nombre y apellido de un autor, y editorial"
$ws.Range("C218").Value = "Consulta que devuelve todos los items prestados en orden alfabetico, filtrable por nombre y apellido de una persona, dias de prestamo, tolerancia y maximo de tolerancia"

$ws.Range("A222").Value = "NotBorrowedTotal"
$ws.Range("A220").Value = "NotBorrowed"
$ws.Range("A224").Value = "TopMostBorrowed"
$ws.Range("A226").Value = "MostBorrowedPerMonth"
$ws.Range("A228").Value = "AgeOfPeopleLoan"

$ws.Range("C220").Value = "Consulta que devuelve un cursor con todos los items que no se encuentran bajo ningun prestamo"
$ws.Range("C222").Value = "Consulta que devuelve la cantidad de todos los items que no se encuentran bajo ningun prestamo"
$ws.Range("C224").Value = "Consulta que devuelve la cantidad de N items que se encuentran prestados"
$ws.Range("C226").Value = "Consulta que devuelve un cursor con todos los items prestados N veces durante N meses"
$ws.Range("C228").Value = "Consulta que devuelve los items prestados en conjunto de la edad de las personas dividido en grupos de edad"

# E column ("Function") for every new row.
foreach ($r in $pairs) {
    $ws.Range("E${r}").Value = "Function"
}

# F column repeats the same name as the A column on each row (literal,
# not read back from A, since Range.Value round-trips unreliably here).
$ws.Range("F212").Value = "allItemsTotal"
$ws.Range("F214").Value = "allLendedItemsTotal"
$ws.Range("F216").Value = "allItems"
$ws.Range("F218").Value = "allLendedItems"
$ws.Range("F220").Value = "NotBorrowed"
$ws.Range("F222").Value = "NotBorrowedTotal"
$ws.Range("F224").Value = "TopMostBorrowed"
$ws.Range("F226").Value = "MostBorrowedPerMonth"
$ws.Range("F228").Value = "AgeOfPeopleLoan"

# ---------------------------------------------------------------
# 4. Row 218/219 get a taller custom height (longer wrapped text).
# ---------------------------------------------------------------
$ws.Rows("218").RowHeight = 24
$ws.Rows("219").RowHeight = 21.75

# ---------------------------------------------------------------
# 5. Move the active selection to match where the author ended up.
# ---------------------------------------------------------------
$ws.Range("L19").Select()
